# Update vm_pu.xlsx results for the 380 kV case (slack bus voltage 1.05 -> 1.02 p.u.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.016124982111738
$ws.Cells.Item(2, 4).Value = 1.043894830642576
$ws.Cells.Item(2, 5).Value = 1.017677374847453
$ws.Cells.Item(2, 6).Value = 1.048137536470974
$ws.Cells.Item(2, 9).Value = 1.036598690647282
$ws.Cells.Item(2, 10).Value = 1.021346577523786
$ws.Cells.Item(2, 11).Value = 1.046667409750793
$ws.Cells.Item(2, 12).Value = 1.020525852793088
$ws.Cells.Item(2, 13).Value = 1.050898220351085
$ws.Cells.Item(2, 14).Value = 1.022797006813959
$ws.Cells.Item(3, 2).Value = 1.019999999999999
$ws.Cells.Item(3, 3).Value = 1.017205078611163
$ws.Cells.Item(3, 4).Value = 1.044443557837969
$ws.Cells.Item(3, 5).Value = 1.018597832551408
$ws.Cells.Item(3, 6).Value = 1.048899337674078
$ws.Cells.Item(3, 9).Value = 1.036685689491114
$ws.Cells.Item(3, 10).Value = 1.022061300645976
$ws.Cells.Item(3, 11).Value = 1.047027950392333
$ws.Cells.Item(3, 12).Value = 1.021251640192673
$ws.Cells.Item(3, 13).Value = 1.051472136450446
$ws.Cells.Item(3, 14).Value = 1.023512744924963
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.017904506731288
$ws.Cells.Item(4, 4).Value = 1.044798477238836
$ws.Cells.Item(4, 5).Value = 1.019194241347825
$ws.Cells.Item(4, 6).Value = 1.049392436457418
$ws.Cells.Item(4, 9).Value = 1.036740602532535
$ws.Cells.Item(4, 10).Value = 1.022523792877596
$ws.Cells.Item(4, 11).Value = 1.047260393243156
$ws.Cells.Item(4, 12).Value = 1.02172146245629
$ws.Cells.Item(4, 13).Value = 1.051842998256935
$ws.Cells.Item(4, 14).Value = 1.023975893948601
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.018198674171667
$ws.Cells.Item(5, 4).Value = 1.044947649021849
$ws.Cells.Item(5, 5).Value = 1.019445165209651
$ws.Cells.Item(5, 6).Value = 1.049599772204728
$ws.Cells.Item(5, 9).Value = 1.036763356650771
$ws.Cells.Item(5, 10).Value = 1.022718228724144
$ws.Cells.Item(5, 11).Value = 1.047357907083703
$ws.Cells.Item(5, 12).Value = 1.021919020633751
$ws.Cells.Item(5, 13).Value = 1.051998787232672
$ws.Cells.Item(5, 14).Value = 1.024170605916353
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.018248073681443
$ws.Cells.Item(6, 4).Value = 1.04497269343432
$ws.Cells.Item(6, 5).Value = 1.019487307766448
$ws.Cells.Item(6, 6).Value = 1.04963458690519
$ws.Cells.Item(6, 9).Value = 1.036767157717844
$ws.Cells.Item(6, 10).Value = 1.022750875589804
$ws.Cells.Item(6, 11).Value = 1.047374268016943
$ws.Cells.Item(6, 12).Value = 1.021952194141809
$ws.Cells.Item(6, 13).Value = 1.052024937727099
$ws.Cells.Item(6, 14).Value = 1.024203299144307
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.017908436908429
$ws.Cells.Item(7, 4).Value = 1.044800470624689
$ws.Cells.Item(7, 5).Value = 1.01919759344474
$ws.Cells.Item(7, 6).Value = 1.049395206744739
$ws.Cells.Item(7, 9).Value = 1.03674090787694
$ws.Cells.Item(7, 10).Value = 1.022526390923155
$ws.Cells.Item(7, 11).Value = 1.047261697035665
$ws.Cells.Item(7, 12).Value = 1.021724102062431
$ws.Cells.Item(7, 13).Value = 1.051845080395564
$ws.Cells.Item(7, 14).Value = 1.023978495683682
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.016489895410684
$ws.Cells.Item(8, 4).Value = 1.04408030419311
$ws.Cells.Item(8, 5).Value = 1.017988279543569
$ws.Cells.Item(8, 6).Value = 1.04839495598871
$ws.Cells.Item(8, 9).Value = 1.036628377842328
$ws.Cells.Item(8, 10).Value = 1.021588117499817
$ws.Cells.Item(8, 11).Value = 1.046789431530405
$ws.Cells.Item(8, 12).Value = 1.020771096598984
$ws.Cells.Item(8, 13).Value = 1.051092280636033
$ws.Cells.Item(8, 14).Value = 1.023038889804462
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.013994322980224
$ws.Cells.Item(9, 4).Value = 1.042810261852656
$ws.Cells.Item(9, 5).Value = 1.015863551461126
$ws.Cells.Item(9, 6).Value = 1.046633713303387
$ws.Cells.Item(9, 9).Value = 1.036419539538403
$ws.Cells.Item(9, 10).Value = 1.019934916795065
$ws.Cells.Item(9, 11).Value = 1.045950782560976
$ws.Cells.Item(9, 12).Value = 1.019093247164828
$ws.Cells.Item(9, 13).Value = 1.049761988321898
$ws.Cells.Item(9, 14).Value = 1.021383341365085
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.012333338517514
$ws.Cells.Item(10, 4).Value = 1.041962996489378
$ws.Cells.Item(10, 5).Value = 1.014451292319811
$ws.Cells.Item(10, 6).Value = 1.045460554614864
$ws.Cells.Item(10, 9).Value = 1.036273265164189
$ws.Cells.Item(10, 10).Value = 1.018832906616808
$ws.Cells.Item(10, 11).Value = 1.045387424397937
$ws.Cells.Item(10, 12).Value = 1.017975690923866
$ws.Cells.Item(10, 13).Value = 1.048872688817748
$ws.Cells.Item(10, 14).Value = 1.020279766205972
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.011614758077194
$ws.Cells.Item(11, 4).Value = 1.041596009808772
$ws.Cells.Item(11, 5).Value = 1.013840777224535
$ws.Cells.Item(11, 6).Value = 1.044952825116729
$ws.Cells.Item(11, 9).Value = 1.036208263596358
$ws.Cells.Item(11, 10).Value = 1.018355756636741
$ws.Cells.Item(11, 11).Value = 1.045142491245955
$ws.Cells.Item(11, 12).Value = 1.017492020661047
$ws.Cells.Item(11, 13).Value = 1.048487053310095
$ws.Cells.Item(11, 14).Value = 1.019801938618203
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.011347940582448
$ws.Cells.Item(12, 4).Value = 1.041459679234662
$ws.Cells.Item(12, 5).Value = 1.013614155783613
$ws.Cells.Item(12, 6).Value = 1.044764272086699
$ws.Cells.Item(12, 9).Value = 1.036183869836435
$ws.Cells.Item(12, 10).Value = 1.018178526244488
$ws.Cells.Item(12, 11).Value = 1.045051364006105
$ws.Cells.Item(12, 12).Value = 1.017312400179648
$ws.Cells.Item(12, 13).Value = 1.048343728062189
$ws.Cells.Item(12, 14).Value = 1.019624456538465
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.011405169535659
$ws.Cells.Item(13, 4).Value = 1.041488923270792
$ws.Cells.Item(13, 5).Value = 1.013662760009546
$ws.Cells.Item(13, 6).Value = 1.044804715497713
$ws.Cells.Item(13, 9).Value = 1.036189113650387
$ws.Cells.Item(13, 10).Value = 1.018216542565012
$ws.Cells.Item(13, 11).Value = 1.045070917790934
$ws.Cells.Item(13, 12).Value = 1.017350927740339
$ws.Cells.Item(13, 13).Value = 1.0483744755481
$ws.Cells.Item(13, 14).Value = 1.019662526846524
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.011592700918511
$ws.Cells.Item(14, 4).Value = 1.041584740983652
$ws.Cells.Item(14, 5).Value = 1.013822041541607
$ws.Cells.Item(14, 6).Value = 1.044937238431367
$ws.Cells.Item(14, 9).Value = 1.036206252281912
$ws.Cells.Item(14, 10).Value = 1.018341106627248
$ws.Cells.Item(14, 11).Value = 1.04513496165438
$ws.Cells.Item(14, 12).Value = 1.01747717243609
$ws.Cells.Item(14, 13).Value = 1.04847520768577
$ws.Cells.Item(14, 14).Value = 1.019787267804016
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.011708257806043
$ws.Cells.Item(15, 4).Value = 1.041643775458241
$ws.Cells.Item(15, 5).Value = 1.013920200179214
$ws.Cells.Item(15, 6).Value = 1.045018895581832
$ws.Cells.Item(15, 9).Value = 1.036216778946386
$ws.Cells.Item(15, 10).Value = 1.018417855240115
$ws.Cells.Item(15, 11).Value = 1.045174401600659
$ws.Cells.Item(15, 12).Value = 1.017554960763877
$ws.Cells.Item(15, 13).Value = 1.048537261121052
$ws.Cells.Item(15, 14).Value = 1.019864125408717
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.012381041717336
$ws.Cells.Item(16, 4).Value = 1.041987349924349
$ws.Cells.Item(16, 5).Value = 1.014491831353991
$ws.Cells.Item(16, 6).Value = 1.045494256500886
$ws.Cells.Item(16, 9).Value = 1.036277544108218
$ws.Cells.Item(16, 10).Value = 1.01886457408961
$ws.Cells.Item(16, 11).Value = 1.045403658951476
$ws.Cells.Item(16, 12).Value = 1.018007795585506
$ws.Cells.Item(16, 13).Value = 1.048898270446034
$ws.Cells.Item(16, 14).Value = 1.020311478650217
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.012803231496654
$ws.Cells.Item(17, 4).Value = 1.042202835632306
$ws.Cells.Item(17, 5).Value = 1.014850669180482
$ws.Cells.Item(17, 6).Value = 1.04579250768644
$ws.Cells.Item(17, 9).Value = 1.036315215525682
$ws.Cells.Item(17, 10).Value = 1.019144796573791
$ws.Cells.Item(17, 11).Value = 1.045547200522158
$ws.Cells.Item(17, 12).Value = 1.018291911085552
$ws.Cells.Item(17, 13).Value = 1.049124572418987
$ws.Cells.Item(17, 14).Value = 1.020592099082468
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.013049549119455
$ws.Cells.Item(18, 4).Value = 1.042328513557489
$ws.Cells.Item(18, 5).Value = 1.015060069903877
$ws.Cells.Item(18, 6).Value = 1.045966497128662
$ws.Cells.Item(18, 9).Value = 1.03633702804794
$ws.Cells.Item(18, 10).Value = 1.019308248267232
$ws.Cells.Item(18, 11).Value = 1.045630829762803
$ws.Cells.Item(18, 12).Value = 1.018457653887013
$ws.Cells.Item(18, 13).Value = 1.049256516111987
$ws.Cells.Item(18, 14).Value = 1.020755782896061
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.013133547501482
$ws.Cells.Item(19, 4).Value = 1.042371364516565
$ws.Cells.Item(19, 5).Value = 1.01513148652601
$ws.Cells.Item(19, 6).Value = 1.046025827143363
$ws.Cells.Item(19, 9).Value = 1.036344438296682
$ws.Cells.Item(19, 10).Value = 1.019363981515951
$ws.Cells.Item(19, 11).Value = 1.045659328825534
$ws.Cells.Item(19, 12).Value = 1.01851417179791
$ws.Cells.Item(19, 13).Value = 1.049301496224614
$ws.Cells.Item(19, 14).Value = 1.020811595292387
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.012757928173771
$ws.Cells.Item(20, 4).Value = 1.042179717204134
$ws.Cells.Item(20, 5).Value = 1.014812159273299
$ws.Cells.Item(20, 6).Value = 1.045760505601839
$ws.Cells.Item(20, 9).Value = 1.036311190343522
$ws.Cells.Item(20, 10).Value = 1.01911473107669
$ws.Cells.Item(20, 11).Value = 1.045531809803189
$ws.Cells.Item(20, 12).Value = 1.018261425796254
$ws.Cells.Item(20, 13).Value = 1.04910029797675
$ws.Cells.Item(20, 14).Value = 1.020561990888912
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.011537474973425
$ws.Cells.Item(21, 4).Value = 1.041556525464579
$ws.Cells.Item(21, 5).Value = 1.013775132927139
$ws.Cells.Item(21, 6).Value = 1.044898212602746
$ws.Cells.Item(21, 9).Value = 1.036201212258058
$ws.Cells.Item(21, 10).Value = 1.018304425493793
$ws.Cells.Item(21, 11).Value = 1.045116106411786
$ws.Cells.Item(21, 12).Value = 1.017439995518422
$ws.Cells.Item(21, 13).Value = 1.04844554686225
$ws.Cells.Item(21, 14).Value = 1.019750534579144
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.010770678651338
$ws.Cells.Item(22, 4).Value = 1.041164612041344
$ws.Cells.Item(22, 5).Value = 1.013123987446477
$ws.Cells.Item(22, 6).Value = 1.044356288919601
$ws.Cells.Item(22, 9).Value = 1.036130622730941
$ws.Cells.Item(22, 10).Value = 1.017794979510281
$ws.Cells.Item(22, 11).Value = 1.044853880651865
$ws.Cells.Item(22, 12).Value = 1.016923739226637
$ws.Cells.Item(22, 13).Value = 1.04803339910091
$ws.Cells.Item(22, 14).Value = 1.019240365123901
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.01117711978751
$ws.Cells.Item(23, 4).Value = 1.041372380467145
$ws.Cells.Item(23, 5).Value = 1.013469088978251
$ws.Cells.Item(23, 6).Value = 1.044643550104049
$ws.Cells.Item(23, 9).Value = 1.036168180026246
$ws.Cells.Item(23, 10).Value = 1.018065044001812
$ws.Cells.Item(23, 11).Value = 1.044992972219733
$ws.Cells.Item(23, 12).Value = 1.017197396497385
$ws.Cells.Item(23, 13).Value = 1.048251931393134
$ws.Cells.Item(23, 14).Value = 1.019510813137988
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.012778398605191
$ws.Cells.Item(24, 4).Value = 1.042190163461667
$ws.Cells.Item(24, 5).Value = 1.014829559943921
$ws.Cells.Item(24, 6).Value = 1.045774965890625
$ws.Cells.Item(24, 9).Value = 1.036313009646311
$ws.Cells.Item(24, 10).Value = 1.019128316372385
$ws.Cells.Item(24, 11).Value = 1.045538764503777
$ws.Cells.Item(24, 12).Value = 1.018275200715029
$ws.Cells.Item(24, 13).Value = 1.049111266720609
$ws.Cells.Item(24, 14).Value = 1.020575595477286
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.014639006152575
$ws.Cells.Item(25, 4).Value = 1.043138707059503
$ws.Cells.Item(25, 5).Value = 1.016412102081107
$ws.Cells.Item(25, 6).Value = 1.047088868053269
$ws.Cells.Item(25, 9).Value = 1.036474774699796
$ws.Cells.Item(25, 10).Value = 1.019934916795065
$ws.Cells.Item(25, 11).Value = 1.045950782560976
$ws.Cells.Item(25, 12).Value = 1.019093247164828
$ws.Cells.Item(25, 13).Value = 1.049761988321898
$ws.Cells.Item(25, 14).Value = 1.021811319835597

Write-Host "case with 380 kV done"
